# Add a new "Save" column (H) to the s_vals sheet, mirroring the header
# style of the existing columns (e.g. G: "sum") and writing a numeric 0
# for the data row, matching the shape of the existing "Win" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (G1) onto the new
# header cell (H1) so it picks up the same bold/centered/bordered style,
# then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data cell for row 2, plain numeric value.
$ws.Range("H2").Value = 0
